# Fruta / hortaliza, semanal
# Insert a new weekly record at row 437 (shifts existing rows 437.. down by one),
# matching the rest of the block's constant columns (A,B,C,E-J) and populating
# the new price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 437; this pushes rows 437-509 down to 438-510
$ws.Rows.Item(437).Insert()

# Populate the newly inserted row 437 with the new weekly observation.
$ws.Cells.Item(437, 1).Value = 7
$ws.Cells.Item(437, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(437, 3).Value = "Ñuble"
$ws.Cells.Item(437, 4).Value = 45275
$ws.Cells.Item(437, 4).NumberFormat = $ws.Cells.Item(438, 4).NumberFormat
$ws.Cells.Item(437, 5).Value = 16
$ws.Cells.Item(437, 6).Value = "Fruta"
$ws.Cells.Item(437, 7).Value = 100102
$ws.Cells.Item(437, 8).Value = "Cítricos"
$ws.Cells.Item(437, 9).Value = 100102004
$ws.Cells.Item(437, 10).Value = "Mandarina"
$ws.Cells.Item(437, 11).Value = "Murcott"
$ws.Cells.Item(437, 12).Value = "Primera"
$ws.Cells.Item(437, 13).Value = 120
$ws.Cells.Item(437, 14).Value = 11000
$ws.Cells.Item(437, 15).Value = 12000
$ws.Cells.Item(437, 16).Value = 11500
$ws.Cells.Item(437, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(437, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(437, 19).Value = 639
$ws.Cells.Item(437, 20).Value = 18
